# Generate Report for Archive
#
# 1) Update the localization status text from "Ready for handoff" to
#    "In Translation" everywhere it appears (Overview!E2:F4 and the
#    Status column (C2:C4) on each language sheet).
# 2) Narrow the status columns (Overview!E:F, zh-cn!C, de-de!C) to the
#    new, narrower width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- 1) Status text: "Ready for handoff" -> "In Translation" ---------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E4").Value = "In Translation"
$overview.Range("F4").Value = "In Translation"

$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"

# --- 2) Column widths: 17.2159881591797 -> 13.4101845877511 ----------------
# ColumnWidth is quantized by the engine to 1/6-character increments, so we
# pick the ColumnWidth input (12.5 "characters") whose rounded, persisted
# width is the closest achievable value to the target 13.4101845877511.
$newColumnWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newColumnWidth  # column E
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth  # column F

$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C

$dede.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C
